# Update MLD ("Base de données/ChampTable.docx")
# - id_droits -> id_droit (singulariser le champ)
# - _GoBack bookmark déplacé de "nom_departement" vers "Impact|es"
# - Impacte -> Impactes (titre)
# - Priorite -> Priorites (titre)

$d = $word.ActiveDocument

# 1. "id_droits" -> "id_droit"
$d.Content.Find.Execute("id_droits", $true, $false, $false, $false, $false, $true, 1, $false, "id_droit", 2) | Out-Null

# 2. Titre "Impacte" -> "Impactes" (match exact, sensible à la casse pour ne
#    toucher que le titre, pas "id_impacte" / "Description_impacte")
$d.Content.Find.Execute("Impacte", $true, $false, $false, $false, $false, $true, 1, $false, "Impactes", 2) | Out-Null

# 3. Titre "Priorite" -> "Priorites"
$d.Content.Find.Execute("Priorite", $true, $false, $false, $false, $false, $true, 1, $false, "Priorites", 2) | Out-Null

# 4. Déplacer le signet caché "_GoBack" : il était après "nom_departement",
#    il se retrouve maintenant au milieu de "Impactes" (entre "Impact" et "es")
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rng = $d.Content
$found = $rng.Find.Execute("Impact", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bmRange = $d.Range($rng.End, $rng.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
